# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F2: 5683 -> 5685
# F3: 5    -> 6
# F6: 37   -> 41

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5685
    $ws.Range("F3").Value = 6
    $ws.Range("F6").Value = 41
}
